$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged), update values
$ws.Range("B3").Value = 0.01506966406258767
$ws.Range("C3").Value = 0.01528788969539829
$ws.Range("D3").Value = 0.01554042569242123

# Row 4 - rename GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.01404973796260023
$ws.Range("C4").Value = 0.01393096764415543
$ws.Range("D4").Value = 0.01396461944426478

# Row 5 - rename AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.01341175649887488
$ws.Range("C5").Value = 0.01351320911982408
$ws.Range("D5").Value = 0.01270016220192608
